$d = $word.ActiveDocument

# 1) Trim the "Processos desejados" line down to just the label (drop "Reembolso")
$d.Content.Find.Execute(
    "Processos desejados: Reembolso", $true, $false, $false, $false, $false,
    $true, 1, $false, "Processos desejados: ", 2) | Out-Null

# 2) Trim the "Informacoes necessarias pelo ERP" line down to just the label
#    (drop "Valor total relatorio, CPF")
$d.Content.Find.Execute(
    "Informações necessárias pelo ERP: Valor total relatório, CPF", $true, $false, $false, $false, $false,
    $true, 1, $false, "Informações necessárias pelo ERP: ", 2) | Out-Null

# 3) Replace the big analysis paragraph body with the new write-up. Build it as an
#    array of lines and join with a manual line break (Chr(11), same as Word's <w:br/>)
#    so each line becomes its own run-text segment separated by a <w:br/>.
$analysisLines = @(
    "Para realizar a integração do ERP SAP ECC/4HANA com o seu SaaS Paytrack seguindo as diretrizes fornecidas, você pode seguir o seguinte modelo de análise funcional:",
    "",
    "### Integração do ERP SAP ECC/4HANA com SaaS Paytrack",
    "",
    "#### Requisitos de Integração:",
    "1. Informações e campos necessários pelo ERP:",
    "   - Lista dos campos obrigatórios e opcionais a serem integrados.",
    "   - Exemplo: bukrs (empresa), kunnr (cliente), belnr (número do documento), etc.",
    "",
    "#### Mapeamento de Campos:",
    "| Campo ERP   | Descrição             | Campo SaaS Paytrack |",
    "|-------------|-----------------------|---------------------|",
    "| bukrs       | Empresa               | Company             |",
    "| kunnr       | Cliente               | Customer            |",
    "| belnr       | Número do Documento   | DocumentNumber      |",
    "| ...         | ...                   | ...                 |",
    "",
    "#### Formato de Saída JSON de Exemplo:",
    "``````json",
    "{",
    "   `"Company`": `"12345`",",
    "   `"Customer`": `"67890`",",
    "   `"DocumentNumber`": `"ABC123`",",
    "   ...",
    "}",
    "``````",
    "",
    "#### Cenários de Integração:",
    "1. Adiantamento:",
    "   - Descrição dos campos específicos necessários para este cenário.",
    "   - Mapeamento de campos relacionados ao adiantamento.",
    "   - Exemplo de JSON para adiantamento.",
    "",
    "2. Prestação de Contas:",
    "   - Descrição dos campos específicos necessários para este cenário.",
    "   - Mapeamento de campos relacionados à prestação de contas.",
    "   - Exemplo de JSON para prestação de contas.",
    "",
    "#### Observações:",
    "1. Utilização de comunicação síncrona com os Webservices do cliente.",
    "2. Paytrack ativa nas integrações, aguardando disponibilização de Webservice pelo cliente.",
    "3. Documentação clara e separada por cenário para facilitar a implementação.",
    "",
    "Ao seguir este modelo de análise funcional, você terá um documento completo e estruturado para guiar o desenvolvimento da integração entre o ERP SAP ECC/4HANA e o seu SaaS Paytrack, considerando os requisitos e diretrizes específicas fornecidas."
)
$analysisPara = $d.Paragraphs.Item(4)
$analysisPara.Range.Text = $analysisLines -join [char]11

# 4) Remove the two data rows from the field-mapping table, keep only the header row
$tbl = $d.Tables.Item(1)
while ($tbl.Rows.Count -gt 1) {
    $tbl.Rows.Item($tbl.Rows.Count).Delete()
}
